$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04642033333333333
$ws.Range("H2").Value = 0.139261
$ws.Range("I2").Value = 0.002653834138691699
$ws.Range("J2").Value = 0.0026538341386917
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 0.1280989678014444
$ws.Range("R2").Value = 1.152890710213
$ws.Range("S2").Value = 0.0006831147774960661
$ws.Range("T2").Value = 0.0006831147774960662
$ws.Range("G3").Value = 0.04642033333333333
$ws.Range("H3").Value = 0.139261
$ws.Range("I3").Value = 0.002653834138691699
$ws.Range("J3").Value = 0.0026538341386917
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 0.3363258524156666
$ws.Range("R3").Value = 3.026932671741
$ws.Range("S3").Value = 0.001793528580146079
$ws.Range("T3").Value = 0.001793528580146079
$ws.Range("G4").Value = 0.04642033333333333
$ws.Range("H4").Value = 0.139261
$ws.Range("I4").Value = 0.002653834138691699
$ws.Range("J4").Value = 0.0026538341386917
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 0.02414002783711111
$ws.Range("R4").Value = 0.217260250534
$ws.Range("S4").Value = 0.0001287317925173091
$ws.Range("T4").Value = 0.0001287317925173091
$ws.Range("G5").Value = 0.04642033333333333
$ws.Range("H5").Value = 0.139261
$ws.Range("I5").Value = 0.002653834138691699
$ws.Range("J5").Value = 0.0026538341386917
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 0.009087120665777779
$ws.Range("R5").Value = 0.081784085992
$ws.Range("S5").Value = [double]"4.845898853224558E-05"
$ws.Range("T5").Value = [double]"4.845898853224559E-05"
$ws.Range("I6").Value = 0.7487035686457026
$ws.Range("J6").Value = 0.7487035686457028
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 36.13946815080677
$ws.Range("R6").Value = 325.255213357261
$ws.Range("S6").Value = 0.192721340135468
$ws.Range("T6").Value = 0.192721340135468
$ws.Range("I7").Value = 0.7487035686457026
$ws.Range("J7").Value = 0.7487035686457028
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.5059929062052914
$ws.Range("T7").Value = 0.5059929062052915
$ws.Range("I8").Value = 0.7487035686457026
$ws.Range("J8").Value = 0.7487035686457028
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 6.810419960066445
$ws.Range("R8").Value = 61.293779640598
$ws.Range("S8").Value = 0.03631800158520167
$ws.Range("T8").Value = 0.03631800158520167
$ws.Range("I9").Value = 0.7487035686457026
$ws.Range("J9").Value = 0.7487035686457028
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 2.563671772847111
$ws.Range("R9").Value = 23.073045955624
$ws.Range("S9").Value = 0.0136713207197416
$ws.Range("T9").Value = 0.01367132071974161
$ws.Range("G10").Value = 4.148506
$ws.Range("H10").Value = 12.445518
$ws.Range("I10").Value = 0.2371686297104146
$ws.Range("J10").Value = 0.2371686297104146
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 11.44798622409933
$ws.Range("R10").Value = 103.031876016894
$ws.Range("S10").Value = 0.06104880231646538
$ws.Range("T10").Value = 0.06104880231646539
$ws.Range("G11").Value = 4.148506
$ws.Range("H11").Value = 12.445518
$ws.Range("I11").Value = 0.2371686297104146
$ws.Range("J11").Value = 0.2371686297104146
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 30.056867680862
$ws.Range("R11").Value = 270.511809127758
$ws.Range("S11").Value = 0.1602845895672332
$ws.Range("T11").Value = 0.1602845895672333
$ws.Range("G12").Value = 4.148506
$ws.Range("H12").Value = 12.445518
$ws.Range("I12").Value = 0.2371686297104146
$ws.Range("J12").Value = 0.2371686297104146
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 2.157353106521334
$ws.Range("R12").Value = 19.416177958692
$ws.Range("S12").Value = 0.0115045406893993
$ws.Range("T12").Value = 0.0115045406893993
$ws.Range("G13").Value = 4.148506
$ws.Range("H13").Value = 12.445518
$ws.Range("I13").Value = 0.2371686297104146
$ws.Range("J13").Value = 0.2371686297104146
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 0.8121004718773334
$ws.Range("R13").Value = 7.308904246896
$ws.Range("S13").Value = 0.004330697137316664
$ws.Range("T13").Value = 0.004330697137316665
$ws.Range("G14").Value = 0.2007003333333333
$ws.Range("H14").Value = 0.602101
$ws.Range("I14").Value = 0.01147396750519105
$ws.Range("J14").Value = 0.01147396750519105
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.553841467548111
$ws.Range("R14").Value = 4.984573207933
$ws.Range("S14").Value = 0.002953476498410602
$ws.Range("T14").Value = 0.002953476498410603
$ws.Range("G15").Value = 0.2007003333333333
$ws.Range("H15").Value = 0.602101
$ws.Range("I15").Value = 0.01147396750519105
$ws.Range("J15").Value = 0.01147396750519105
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 1.454119473975667
$ws.Range("R15").Value = 13.087075265781
$ws.Range("S15").Value = 0.007754398946112223
$ws.Range("T15").Value = 0.007754398946112224
$ws.Range("G16").Value = 0.2007003333333333
$ws.Range("H16").Value = 0.602101
$ws.Range("I16").Value = 0.01147396750519105
$ws.Range("J16").Value = 0.01147396750519105
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.1043704619437778
$ws.Range("R16").Value = 0.9393341574940001
$ws.Range("S16").Value = 0.0005565775127743182
$ws.Range("T16").Value = 0.0005565775127743182
$ws.Range("G17").Value = 0.2007003333333333
$ws.Range("H17").Value = 0.602101
$ws.Range("I17").Value = 0.01147396750519105
$ws.Range("J17").Value = 0.01147396750519105
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.03928856205244444
$ws.Range("R17").Value = 0.353597058472
$ws.Range("S17").Value = 0.0002095145478939085
$ws.Range("T17").Value = 0.0002095145478939086
